# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new Price text (kept as literal text, matching existing inline-string cells)
$priceUpdates = @{
    2 = "29.447.22"
    3 = "1.883.25"
    4 = "1.001"
    5 = "0.7204"
    6 = "243.61"
    7 = "1.001"
    8 = "0.07976"
    9 = "0.3154"
    10 = "25.06"
    11 = "0.08144"
    12 = "1.894.49"
    13 = "5.257"
    14 = "94.93"
    15 = "0.7120"
    16 = "6.412"
    17 = "0.000008468"
    18 = "29.452.65"
    19 = "254.91"
    20 = "13.36"
    21 = "2.138.71"
    22 = "1.000"
    23 = "7.798"
    26 = "9.091"
    27 = "162.61"
    28 = "19.07"
    29 = "1.509"
    30 = "4.427"
    31 = "4.291"
    33 = "0.05326"
    34 = "1.953"
    35 = "0.7575"
    36 = "1.183"
    37 = "2.700"
    38 = "0.01893"
    39 = "1.272.99"
    40 = "2.768"
    41 = "6.477"
    42 = "113.19"
    43 = "74.58"
    44 = "0.9047"
    47 = "2.038.70"
    48 = "1.808"
    49 = "0.5195"
    50 = "9.540"
    51 = "0.4377"
}

# Map of row -> new Volume(1h) text
$volumeUpdates = @{
    2 = "  +0.55%  "
    3 = "  +0.59%  "
    4 = "  +0.07%  "
    5 = "  +1.76%  "
    6 = "  +0.82%  "
    7 = "  +0.13%  "
    8 = "  +2.42%  "
    9 = "  +1.48%  "
    10 = "  +0.20%  "
    11 = "  -2.93%  "
    12 = "  +0.82%  "
    13 = "  +0.42%  "
    14 = "  +4.27%  "
    15 = "  -0.65%  "
    16 = "  +4.81%  "
    17 = "  +1.65%  "
    18 = "  +0.52%  "
    19 = "  +6.11%  "
    20 = "  +1.22%  "
    21 = "  +0.44%  "
    22 = "  +0.04%  "
    23 = "  +0.58%  "
    24 = "  +0.08%  "
    25 = "  -0.24%  "
    26 = "  +0.74%  "
    27 = "  +0.05%  "
    28 = "  +3.14%  "
    29 = "  +0.28%  "
    30 = "  +0.54%  "
    31 = "  -0.50%  "
    32 = "  -2.35%  "
    33 = "  -0.67%  "
    34 = "  +0.87%  "
    35 = "  +1.18%  "
    36 = "  +0.71%  "
    37 = "  +0.74%  "
    38 = "  +0.94%  "
    39 = "  +2.67%  "
    40 = "  +1.39%  "
    41 = "  -0.66%  "
    42 = "  +3.81%  "
    43 = "  +3.11%  "
    44 = "  +1.46%  "
    45 = "  +4.58%  "
    46 = "  +0.13%  "
    47 = "  +0.73%  "
    48 = "  +1.10%  "
    49 = "  -0.07%  "
    50 = "  +1.14%  "
    51 = "  +1.03%  "
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    # Force text storage so numeric-looking strings (e.g. "1.001") are not
    # auto-converted to numbers, matching the original inline-string cells.
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.Style = "Normal"
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeUpdates[$row]
}
